$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.463
$ws.Range("A3").Value = -21.682
$ws.Range("C5").Value = -12.463
$ws.Range("A14").Value = -21.04
$ws.Range("A16").Value = -21.395
$ws.Range("C16").Value = -12.016
$ws.Range("A21").Value = -21.04
$ws.Range("A23").Value = -21.701
$ws.Range("A25").Value = -22.05
